# Reorder columns (weights changed) and replace row 3 with the closest
# counterfactual values, then drop the now-superfluous counterfactual rows
# 4 and 5 ("pesi e closest counterfactual").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: headers in their new column order -----------------------------
$ws.Cells.Item(1,1).Value  = "ER Triage"
$ws.Cells.Item(1,2).Value  = "Release D"
$ws.Cells.Item(1,3).Value  = "IV Liquid"
$ws.Cells.Item(1,4).Value  = "Release B"
$ws.Cells.Item(1,5).Value  = "Leucocytes"
$ws.Cells.Item(1,6).Value  = "Release C"
$ws.Cells.Item(1,7).Value  = "ER Sepsis Triage"
$ws.Cells.Item(1,8).Value  = "other"
$ws.Cells.Item(1,9).Value  = "Admission IC"
$ws.Cells.Item(1,10).Value = "IV Antibiotics"
$ws.Cells.Item(1,11).Value = "LacticAcid"
$ws.Cells.Item(1,12).Value = "ER Registration"
$ws.Cells.Item(1,13).Value = "Admission NC"
$ws.Cells.Item(1,14).Value = "Release A"
$ws.Cells.Item(1,15).Value = "CRP"
# P1 "Label" and Q1 "Type" stay as they are

# --- Row 2: "Original" record, reordered to match the new columns ---------
$ws.Cells.Item(2,1).Value  = 1
$ws.Cells.Item(2,2).Value  = 0
$ws.Cells.Item(2,3).Value  = 0
$ws.Cells.Item(2,4).Value  = 0
$ws.Cells.Item(2,5).Value  = 1
$ws.Cells.Item(2,6).Value  = 0
$ws.Cells.Item(2,7).Value  = 1
$ws.Cells.Item(2,8).Value  = 0
$ws.Cells.Item(2,9).Value  = 0
$ws.Cells.Item(2,10).Value = 1
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 0
$ws.Cells.Item(2,14).Value = 0
$ws.Cells.Item(2,15).Value = 1
# P2 "deviant" and Q2 "Original" stay as they are

# --- Row 3: now holds the single closest counterfactual --------------------
$ws.Cells.Item(3,1).Value  = 1
$ws.Cells.Item(3,2).Value  = 0
$ws.Cells.Item(3,3).Value  = 0
$ws.Cells.Item(3,4).Value  = 0
$ws.Cells.Item(3,5).Value  = 1
$ws.Cells.Item(3,6).Value  = 0
$ws.Cells.Item(3,7).Value  = 1
$ws.Cells.Item(3,8).Value  = 0
$ws.Cells.Item(3,9).Value  = 0
$ws.Cells.Item(3,10).Value = 1
$ws.Cells.Item(3,11).Value = 1
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 13
$ws.Cells.Item(3,14).Value = 0
$ws.Cells.Item(3,15).Value = 1
# P3 "regular" and Q3 "Counterfactual" stay as they are

# --- Drop the other two counterfactual rows --------------------------------
$ws.Rows("4:5").Delete()
